$wb = $excel.ActiveWorkbook
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Rows("4:4").Delete()
$n = $ws4.UsedRange.Rows.Count
for ($r = 4; $r -le $n; $r++) {
    $ws4.Cells.Item($r, 1).Value2 = $r - 1
}
Write-Host "done"
